# Update "想去人数" (F column) counts across sheets to reflect refreshed
# scrape values (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 625
$ws1.Range("F4").Value = 907
$ws1.Range("F5").Value = 666
$ws1.Range("F6").Value = 803
$ws1.Range("F7").Value = 372
$ws1.Range("F8").Value = 574
$ws1.Range("F9").Value = 114
$ws1.Range("F10").Value = 1165
$ws1.Range("F11").Value = 592
$ws1.Range("F12").Value = 353
$ws1.Range("F13").Value = 478
$ws1.Range("F14").Value = 151
$ws1.Range("F15").Value = 52
$ws1.Range("F16").Value = 314
$ws1.Range("F18").Value = 73
$ws1.Range("F20").Value = 39
$ws1.Range("F21").Value = 539
$ws1.Range("F22").Value = 17
$ws1.Range("F23").Value = 548
$ws1.Range("F24").Value = 1

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 74
$ws2.Range("F3").Value = 57
$ws2.Range("F4").Value = 307
$ws2.Range("F9").Value = 205
$ws2.Range("F13").Value = 39

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 74
$ws4.Range("F5").Value = 57
$ws4.Range("F6").Value = 307
$ws4.Range("F7").Value = 625
$ws4.Range("F8").Value = 907
$ws4.Range("F9").Value = 666
$ws4.Range("F10").Value = 803
$ws4.Range("F11").Value = 372
$ws4.Range("F12").Value = 574
$ws4.Range("F13").Value = 114
$ws4.Range("F14").Value = 1165
$ws4.Range("F15").Value = 592
$ws4.Range("F18").Value = 353
$ws4.Range("F19").Value = 478
$ws4.Range("F21").Value = 151
$ws4.Range("F22").Value = 52
$ws4.Range("F24").Value = 314
$ws4.Range("F26").Value = 73
$ws4.Range("F27").Value = 205
$ws4.Range("F32").Value = 39
$ws4.Range("F33").Value = 39
$ws4.Range("F34").Value = 539
$ws4.Range("F35").Value = 17
$ws4.Range("F36").Value = 548
$ws4.Range("F37").Value = 1
